$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (D = Fecha serial, M = Volumen, Q = Unidad de comercialización,
# S = Precio $/Kg, T = Kg / unidad)
$rows = @(
    @{ Row = 2;  D = 44309; M = 80;  Q = "`$/caja 14 kilos granel";    S = 821;   T = 14 },
    @{ Row = 3;  D = 44313; M = 120; Q = "`$/caja 10 kilos empedrada"; S = 11500; T = 1 },
    @{ Row = 4;  D = 44302; M = 80;  Q = "`$/caja 10 kilos empedrada"; S = 11500; T = 1 },
    @{ Row = 5;  D = 44323; M = 80;  Q = "`$/caja 10 kilos empedrada"; S = 11500; T = 1 },
    @{ Row = 6;  D = 44330; M = 60;  Q = "`$/caja 10 kilos empedrada"; S = 11500; T = 1 },
    @{ Row = 7;  D = 44322; M = 60;  Q = "`$/caja 10 kilos empedrada"; S = 11500; T = 1 },
    @{ Row = 8;  D = 44306; M = 80;  Q = "`$/caja 10 kilos empedrada"; S = 11500; T = 1 },
    @{ Row = 9;  D = 44316; M = 120; Q = "`$/caja 10 kilos empedrada"; S = 11500; T = 1 },
    @{ Row = 10; D = 44327; M = 60;  Q = "`$/caja 10 kilos empedrada"; S = 11500; T = 1 }
)

foreach ($r in $rows) {
    $ws.Range("D$($r.Row)").Value = $r.D
    $ws.Range("M$($r.Row)").Value = $r.M
    $ws.Range("Q$($r.Row)").Value = $r.Q
    $ws.Range("S$($r.Row)").Value = $r.S
    $ws.Range("T$($r.Row)").Value = $r.T
}
